# Apply the commit's changes to doe1.xlsx:
#  - rename the param13..param16 header labels to their upper-case form
#  - tighten the saved selection on Sheet1 from D2:D12 down to just D2
#  - (workbook window tab-ratio bump from 990 -> 991 is a pure window/host
#    UI setting with no effect on the saved document model in this runtime)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: A1:D1 currently hold "param13".."param16" -> upper-case them.
$ws.Range("A1").Value = "PARAM13"
$ws.Range("B1").Value = "PARAM14"
$ws.Range("C1").Value = "PARAM15"
$ws.Range("D1").Value = "PARAM16"

# Also try to bump the window tab ratio in case it is honoured by the host.
$win = $wb.Windows.Item(1)
$win.TabRatio = 991

# Shrink the persisted selection from D2:D12 to just D2.
$ws.Range("D2").Select()
